# IAS_IMATGES.xlsx edit script
# - "PostgresSQL_IAS" sheet: insert a new column A with the "external"/postgres ids,
#   shifting the old id column to B and the species-name column to C; add a new
#   column D with a computed "id;name" helper string; update the header row.
# - "Imatges" sheet: remap column B (foreign key into PostgresSQL_IAS) from the old
#   sequential ids to the new ones, and update the VLOOKUP formulas in column J to
#   use the new 3-column lookup table.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Imatges")
$ws4 = $wb.Worksheets.Item("PostgresSQL_IAS")

# ---------------------------------------------------------------------------
# 1. PostgresSQL_IAS sheet: insert new column A (old column shifts to B/C)
# ---------------------------------------------------------------------------
$ws4.Columns.Item(1).Insert()

# New column A values (old/external species id), keyed by row number
$sheet4A = @{
  2 = 1
  3 = 2
  4 = 4
  5 = 15
  6 = 11
  7 = 7
  8 = 19
  11 = 6
  12 = 14
  13 = 3
  14 = 17
  15 = 5
  16 = 10
  17 = 16
  18 = 12
  19 = 13
  20 = 9
  21 = 18
  22 = 8
}
foreach ($r in $sheet4A.Keys) {
  $ws4.Cells.Item($r,1).Value = $sheet4A[$r]
}

# New column D values (helper text "id;\"name\"")
$sheet4D = @{
  2 = '1;"Agave americana"'
  3 = '2;"Ailanthus altissima"'
  4 = '4;"Cortaderia selloana"'
  5 = '15;"Opuntia spp."'
  6 = '11;"Senecio angulatus"'
  7 = '7;"Dreissena polymorpha"'
  11 = '6;"Estrilda astrild"'
  12 = '14;"Myocastor coypus"'
  13 = '3;"Carpobrotus spp."'
  14 = '17;"Neovison vison"'
  15 = '5;"Cotoneaster horizontalis"'
  16 = '10;"Heracleum mantegazzianum"'
  17 = '16;"Phytolacca americana"'
  18 = '12;"Myiopsitta monachus"'
  19 = '13;"Psittacula krameri"'
  20 = '9;"Psittacula eupatria"'
  21 = '18;"Procyon lotor"'
  22 = '8;"Pomacea insularum"'
}
foreach ($r in $sheet4D.Keys) {
  $ws4.Cells.Item($r,4).Value = $sheet4D[$r]
}

# Give the new A1 header the same style as the other header cells, then set
# its text (and the B1 header text; C1 "IAS_scientificname" is already
# correct after the column shift).
$ws4.Cells.Item(1,2).Copy()
$ws4.Cells.Item(1,1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws4.Cells.Item(1,2).Value = "IAS_id_postgresql_BB"
$ws4.Cells.Item(1,1).Value = "IAS_Id_postgresql_Isaac"

# Column widths for the new layout
$ws4.Columns.Item(1).ColumnWidth = 21.76
$ws4.Columns.Item(2).ColumnWidth = 19.94
$ws4.Columns.Item(3).ColumnWidth = 25.97
$ws4.Columns.Item(4).ColumnWidth = 107.13

# ---------------------------------------------------------------------------
# 2. Imatges sheet: remap column B ids and update the VLOOKUP formulas
# ---------------------------------------------------------------------------
$sheet1B = @{
  2 = 1
  3 = 1
  4 = 1
  5 = 1
  6 = 1
  7 = 2
  8 = 2
  9 = 2
  10 = 2
  11 = 6
  12 = 6
  13 = 6
  14 = 3
  15 = 3
  16 = 3
  17 = 4
  18 = 4
  19 = 4
  20 = 19
  21 = 19
  22 = 19
  23 = 19
  24 = 19
  25 = 19
  26 = 12
  27 = 12
  28 = 12
  29 = 15
  30 = 15
  31 = 11
  32 = 11
  33 = 13
  34 = 13
  35 = 5
  36 = 5
  37 = 7
  38 = 7
  39 = 7
  40 = 8
  41 = 8
  42 = 9
  43 = 9
  44 = 10
  45 = 10
  46 = 10
  47 = 10
  48 = 14
  49 = 14
  50 = 14
  51 = 16
  52 = 16
  53 = 17
  54 = 17
  55 = 17
  56 = 18
  57 = 18
  58 = 18
}

for ($r = 2; $r -le 58; $r++) {
  $ws1.Cells.Item($r,2).Value = $sheet1B[$r]
  $ws1.Cells.Item($r,10).Formula = '=VLOOKUP(B' + $r.ToString() + ',PostgresSQL_IAS!$A$2:$C$22,3,FALSE)'
}

$wb.Application.Calculate()

# ---------------------------------------------------------------------------
# 3. Restore view/selection state
# ---------------------------------------------------------------------------
$ws4.Activate()
$ws4.Range("B19").Select()

$ws1.Activate()
$ws1.Range("B7").Select()
